$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new program rows (GeneXpert cost lines)
$ws.Range("A8").Value = "econ_program_totalcost_xpert"
$ws.Range("B8").Value = "yes"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "yes"

$ws.Range("A9").Value = "econ_program_unitcost_xpert"
$ws.Range("B9").Value = "yes"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "yes"

# Update the active selection to match the committed state
$ws.Range("E8").Select()
